$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete column D ("BL HOUSE") - merchandise document no longer tracks a
# separate BL HOUSE reference, so the whole column is removed and the
# remaining columns shift left.
$ws.Range("D1").EntireColumn.Delete()

# Refresh the autofilter range (and its backing defined name) to the new
# used range A1:N1 now that a column has been removed.
$ws.AutoFilterMode = $false
$null = $ws.Range("A1:N1").AutoFilter()

foreach ($n in $wb.Names) {
    if ($n.Name -like "*_FilterDatabase*") {
        $n.RefersTo = "=Hoja1!`$A`$1:`$N`$1"
    }
}

# Leave the selection where the editor left it after the rearrangement.
$null = $ws.Range("L8").Select()
